# Auto-generated edit script: apply scheduled-runner market-data refresh
# to the Leve profit sheets (columns H-N) per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1916
$ws.Range("I2").Value = 2547
$ws.Range("J2").Value = 811.75
$ws.Range("K2").Value = 2547
$ws.Range("L2").Value = 811.75
$ws.Range("M2").Value = -2434
$ws.Range("N2").Value = -1037.75
$ws.Range("H33").Value = 185.85185
$ws.Range("I33").Value = 104.26316
$ws.Range("K33").Value = 104.26316
$ws.Range("M33").Value = 124.73684
$ws.Range("H43").Value = 3164.3076
$ws.Range("I43").Value = 2271.1428
$ws.Range("J43").Value = 4206.3335
$ws.Range("K43").Value = 2271.1428
$ws.Range("L43").Value = 4206.3335
$ws.Range("M43").Value = -2202.1428
$ws.Range("N43").Value = -4344.3335
$ws.Range("H58").Value = 2293.9092
$ws.Range("I58").Value = 554.125
$ws.Range("J58").Value = 6933.3335
$ws.Range("K58").Value = 1662.375
$ws.Range("L58").Value = 20800.0005
$ws.Range("M58").Value = -1512.375
$ws.Range("N58").Value = -21100.0005
$ws.Range("H64").Value = 21899.834
$ws.Range("I64").Value = 27266.555
$ws.Range("J64").Value = 5799.6665
$ws.Range("K64").Value = 27266.555
$ws.Range("L64").Value = 5799.6665
$ws.Range("M64").Value = -27018.555
$ws.Range("N64").Value = -6295.6665
$ws.Range("H67").Value = 21899.834
$ws.Range("I67").Value = 27266.555
$ws.Range("J67").Value = 5799.6665
$ws.Range("K67").Value = 27266.555
$ws.Range("L67").Value = 5799.6665
$ws.Range("M67").Value = -26408.555
$ws.Range("N67").Value = -7515.6665
$ws.Range("H70").Value = 6964.9
$ws.Range("J70").Value = 11979.8
$ws.Range("L70").Value = 35939.39999999999
$ws.Range("N70").Value = -36479.39999999999
$ws.Range("H73").Value = 6964.9
$ws.Range("J73").Value = 11979.8
$ws.Range("L73").Value = 35939.39999999999
$ws.Range("N73").Value = -37811.39999999999
$ws.Range("H125").Value = 40214.25
$ws.Range("J125").Value = 3097.4285
$ws.Range("L125").Value = 27876.8565
$ws.Range("N125").Value = -32796.8565
$ws.Range("H132").Value = 18837.375
$ws.Range("I132").Value = 21167.572
$ws.Range("K132").Value = 63502.716
$ws.Range("M132").Value = -60972.716
$ws.Range("H135").Value = 3082.3547
$ws.Range("I135").Value = 2857.52
$ws.Range("J135").Value = 4019.1667
$ws.Range("K135").Value = 25717.68
$ws.Range("L135").Value = 36172.5003
$ws.Range("M135").Value = -23182.68
$ws.Range("N135").Value = -41242.5003

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3749.2693
$ws.Range("I2").Value = 3685.652
$ws.Range("K2").Value = 3685.652
$ws.Range("M2").Value = -3572.652
$ws.Range("H32").Value = 20282.98
$ws.Range("I32").Value = 21824.55
$ws.Range("K32").Value = 21824.55
$ws.Range("M32").Value = -21537.55
$ws.Range("H45").Value = 3498.9412
$ws.Range("I45").Value = 1726.8889
$ws.Range("J45").Value = 5492.5
$ws.Range("K45").Value = 1726.8889
$ws.Range("L45").Value = 5492.5
$ws.Range("M45").Value = -1349.8889
$ws.Range("N45").Value = -6246.5
$ws.Range("H61").Value = 4821.9287
$ws.Range("I61").Value = 1083.3334
$ws.Range("K61").Value = 1083.3334
$ws.Range("M61").Value = -871.3334
$ws.Range("H74").Value = 278761.6
$ws.Range("I74").Value = 375830.88
$ws.Range("K74").Value = 375830.88
$ws.Range("M74").Value = -374956.88
$ws.Range("H77").Value = 278761.6
$ws.Range("I77").Value = 375830.88
$ws.Range("K77").Value = 1879154.4
$ws.Range("M77").Value = -1874786.4
$ws.Range("H116").Value = 3749.2693
$ws.Range("I116").Value = 3685.652
$ws.Range("K116").Value = 3685.652
$ws.Range("M116").Value = -1391.652
$ws.Range("H122").Value = 1935.5
$ws.Range("I122").Value = 1740.3684
$ws.Range("K122").Value = 5221.1052
$ws.Range("M122").Value = -2771.1052
$ws.Range("H136").Value = 4821.9287
$ws.Range("I136").Value = 1083.3334
$ws.Range("K136").Value = 3250.0002
$ws.Range("M136").Value = -700.0001999999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3749.2693
$ws.Range("I3").Value = 3685.652
$ws.Range("K3").Value = 3685.652
$ws.Range("M3").Value = -3571.652
$ws.Range("H86").Value = 1742.7142
$ws.Range("H89").Value = 1742.7142

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2276.923
$ws.Range("I6").Value = 2276.923
$ws.Range("K6").Value = 2276.923
$ws.Range("M6").Value = -2163.923
$ws.Range("H22").Value = 1112.5
$ws.Range("I22").Value = 983.3333
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 983.3333
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -633.3333
$ws.Range("N22").Value = -2200
$ws.Range("H50").Value = 36449.4
$ws.Range("J50").Value = 36449.4
$ws.Range("L50").Value = 36449.4
$ws.Range("N50").Value = -37699.4
$ws.Range("H58").Value = 1153.4231
$ws.Range("J58").Value = 1889.8334
$ws.Range("L58").Value = 1889.8334
$ws.Range("N58").Value = -2295.8334
$ws.Range("H86").Value = 29203.719
$ws.Range("I86").Value = 42112.445
$ws.Range("J86").Value = 12606.786
$ws.Range("K86").Value = 42112.445
$ws.Range("L86").Value = 12606.786
$ws.Range("M86").Value = -40989.445
$ws.Range("N86").Value = -14852.786
$ws.Range("H89").Value = 29203.719
$ws.Range("I89").Value = 42112.445
$ws.Range("J89").Value = 12606.786
$ws.Range("K89").Value = 210562.225
$ws.Range("L89").Value = 63033.93
$ws.Range("M89").Value = -204946.225
$ws.Range("N89").Value = -74265.92999999999
$ws.Range("H107").Value = 1416.5555
$ws.Range("J107").Value = 1241.5
$ws.Range("L107").Value = 1241.5
$ws.Range("N107").Value = -5081.5
$ws.Range("H132").Value = 47084.547
$ws.Range("I132").Value = 63289.188
$ws.Range("K132").Value = 189867.564
$ws.Range("M132").Value = -187337.564
$ws.Range("H136").Value = 1153.4231
$ws.Range("J136").Value = 1889.8334
$ws.Range("L136").Value = 5669.5002
$ws.Range("N136").Value = -10769.5002

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4999.5
$ws.Range("I62").Value = 4999
$ws.Range("K62").Value = 14997
$ws.Range("M62").Value = -14311
$ws.Range("H65").Value = 4999.5
$ws.Range("I65").Value = 4999
$ws.Range("K65").Value = 44991
$ws.Range("M65").Value = -41559
$ws.Range("H131").Value = 124020.74
$ws.Range("J131").Value = 1949.0385
$ws.Range("L131").Value = 5847.1155
$ws.Range("N131").Value = -15927.1155

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 34714.145
$ws.Range("I69").Value = 33999
$ws.Range("K69").Value = 33999
$ws.Range("M69").Value = -33250
$ws.Range("H72").Value = 34714.145
$ws.Range("I72").Value = 33999
$ws.Range("K72").Value = 101997
$ws.Range("M72").Value = -98253
$ws.Range("H99").Value = 8893.666999999999
$ws.Range("J99").Value = 39000
$ws.Range("L99").Value = 39000
$ws.Range("N99").Value = -43492
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 3382.1892
$ws.Range("I122").Value = 3079.1155
$ws.Range("K122").Value = 9237.3465
$ws.Range("M122").Value = -6787.3465
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3530
$ws.Range("N126").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4644.6
$ws.Range("I68").Value = 3979.1667
$ws.Range("J68").Value = 5088.222
$ws.Range("K68").Value = 3979.1667
$ws.Range("L68").Value = 5088.222
$ws.Range("M68").Value = -3230.1667
$ws.Range("N68").Value = -6586.222
$ws.Range("H71").Value = 4644.6
$ws.Range("I71").Value = 3979.1667
$ws.Range("J71").Value = 5088.222
$ws.Range("K71").Value = 19895.8335
$ws.Range("L71").Value = 25441.11
$ws.Range("M71").Value = -16151.8335
$ws.Range("N71").Value = -32929.11

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 65599.8
$ws.Range("J86").Value = 65599.8
$ws.Range("L86").Value = 65599.8
$ws.Range("N86").Value = -67845.8
$ws.Range("H89").Value = 65599.8
$ws.Range("J89").Value = 65599.8
$ws.Range("L89").Value = 327999
$ws.Range("N89").Value = -339231
$ws.Range("H107").Value = 1280.0834
$ws.Range("I107").Value = 1410.25
$ws.Range("J107").Value = 1019.75
$ws.Range("K107").Value = 4230.75
$ws.Range("L107").Value = 3059.25
$ws.Range("M107").Value = -2310.75
$ws.Range("N107").Value = -6899.25
$ws.Range("H122").Value = 44479.97
$ws.Range("I122").Value = 50352.715
$ws.Range("K122").Value = 151058.145
$ws.Range("M122").Value = -148608.145
$ws.Range("H132").Value = 19517.36
$ws.Range("I132").Value = 21539.03
$ws.Range("K132").Value = 64617.09
$ws.Range("M132").Value = -62087.09
$ws.Range("H136").Value = 20888.234
$ws.Range("I136").Value = 22940
$ws.Range("K136").Value = 68820
$ws.Range("M136").Value = -66270
